# Apply updated crypto price/volume figures (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.814.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.59%  '
$ws.Range("D3").Value = "'2.275.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.92%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'230.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").Value = "'0.626"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("D7").Value = "'61.18"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.89%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +4.31%  '
$ws.Range("D10").Value = "'57.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.21%  '
$ws.Range("D11").Value = "'0.0933"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.69%  '
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("D13").Value = "'2.609.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").Value = "'15.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.74%  '
$ws.Range("D15").Value = "'23.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.17%  '
$ws.Range("D16").Value = "'5.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.19%  '
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").Value = "'2.270.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").Value = "'43.720.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.67%  '
$ws.Range("D20").Value = "'0.0₃0935"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.57%  '
$ws.Range("D21").Value = "'73.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("E22").Value = '  +2.61%  '
$ws.Range("D23").Value = "'253.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").Value = "'2.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.38%  '
$ws.Range("D26").Value = "'2.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.91%  '
$ws.Range("E27").Value = '  +1.76%  '
$ws.Range("D28").Value = "'170.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.83%  '
$ws.Range("E29").Value = '  -1.59%  '
$ws.Range("D30").Value = "'20.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = '  +2.37%  '
$ws.Range("D32").Value = "'2.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.09%  '
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("D34").Value = "'5.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("D35").Value = "'4.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.69%  '
$ws.Range("D36").Value = "'0.0658"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.98%  '
$ws.Range("D37").Value = "'6.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'2.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.46%  '
$ws.Range("D39").Value = "'3.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.05%  '
$ws.Range("E40").Value = '  +4.10%  '
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("D42").Value = "'8.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.67%  '
$ws.Range("E43").Value = '  -10.97%  '
$ws.Range("E44").Value = '  +0.91%  '
$ws.Range("D45").Value = "'4.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.55%  '
$ws.Range("E46").Value = '  -1.56%  '
$ws.Range("D47").Value = "'98.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.86%  '
$ws.Range("D48").Value = "'1.478.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("D49").Value = "'16.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.25%  '
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("D51").Value = "'2.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.21%  '
